$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sorted data (by total_registros descending) for rows 2-16.
# Rows 17 and 18 are unchanged.
$data = @(
    @("RUIZ CHIROQUE CLAUDIA JUDITH", 149),
    @("TEMOCHE ECHE URSULA YESSENIA", 147),
    @("BANCAYAN FIESTA DILVER HUMBERTO", 147),
    @("FABIANA REBECA ARRUNATEGUI SILUPU", 147),
    @("GONZALES FIESTAS MARIA MARIBEL", 146),
    @("ANTON INGA FATIMA DEL ROSARIO", 144),
    @("LLENQUE ANTON HELEN JOHANA", 142),
    @("BAUTISTA CHAVESTA ERICKA MEDALIT", 136),
    @("PINTADO CHASQUERO ESTEFANY", 117),
    @("MONDRAGON NONAJULCA MARISOL", 111),
    @("VELASCO PEÑA KAREN ARELLYS", 109),
    @("ORDINOLA JIBAJA JOSE ALBERTO", 103),
    @("FLORES SILUPU MARY CARMEN", 101),
    @("HERNANDEZ CARNERO ARTURO SEBASTIAN", 95),
    @("CASTRO ESTRADA CINTHIA PATRICIA", 84)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
